$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded, which pushes all existing
# records for this product (rows 55-103) down by one row (to 56-104).
# Insert a new row at position 55, shifting everything below it down.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new observation. Most columns
# match the pattern used by the surrounding rows for this market/product.
$ws.Range("A55").Value2 = 7
$ws.Range("B55").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C55").Value2 = "Ñuble"
$ws.Range("D55").Value2 = 45264
$ws.Range("E55").Value2 = 16
$ws.Range("F55").Value2 = 100112026
$ws.Range("G55").Value2 = "Haba"
$ws.Range("H55").Value2 = "Sin especificar"
$ws.Range("I55").Value2 = "Primera"
$ws.Range("J55").Value2 = 60
$ws.Range("K55").Value2 = 10000
$ws.Range("L55").Value2 = 10000
$ws.Range("M55").Value2 = 10000
$ws.Range("N55").Value2 = "$/saco 25 kilos"
$ws.Range("O55").Value2 = "Provincia de Diguillín"
$ws.Range("P55").Value2 = 400
$ws.Range("Q55").Value2 = 25
$ws.Range("R55").Value2 = "Hortaliza"
